$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experimentos")

# Insert a new column before F (shifts old F:I -> G:J), this is the
# "Tolerance" column being added to the results table.
$ws.Columns("F").Insert()

# Row 4 of the table (A6:J6) was missing its Classificador/Configuracoes
# labels - fill them in to match the other rows.
$ws.Range("B6").Value = "SVM"
$ws.Range("C6").Value = "(CODE with parallel)"

# New experiment row (row 5 of the table / row 7 of the sheet).
# C7 is set before F2 so that the new shared string
# "(CODE with parallel - both)" is registered ahead of "Tolerance",
# matching the order they appear in the original authoring.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "SVM"
$ws.Range("C7").Value = "(CODE with parallel - both)"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 200
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 77
$ws.Range("J7").Value = 0.135854341736695

# Header label for the newly inserted column
$ws.Range("F2").Value = "Tolerance"

# Column width adjustments (C got wider to fit the new longer label,
# F takes the same width as the neighbouring columns).
$ws.Columns(3).ColumnWidth = 24.166666666666668
$ws.Columns(6).ColumnWidth = 20.307291666666668

# Move the active selection down to the row below the table
[void]$ws.Range("A8").Select()
